$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 was previously blank (data started at row 2). Fill it in with
# header labels for the card columns - no shifting of existing rows.
$ws.Range("A1").Value = "back"
$ws.Range("B1").Value = "front"

# Update selection to match the edited workbook's view state.
$ws.Range("B1").Select()
